# CIERRE 26 DIC 23
# Update the "VALES DE INSENTIVOS" sheet to reflect the new incentive
# voucher: amount, concept/description and beneficiary name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Voucher amount: 5,000 -> 6,000
$ws.Range("D1").Value = 6000

# Amount spelled out in words: CINCO MIL -> SEIS MIL
$ws.Range("A2").Value = "SEIS  MIL   PESOS 00/100 M.N."

# Concept / reason for the voucher
$ws.Range("A4").Value = "INCENTIVO DEL MES DE  NOVIEMBRRE 2023   "

# Beneficiary name
$ws.Range("C8").Value = "PABLO BAEZ"
